$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Delete the old "couponSortTest" block (rows 9-12) entirely; everything below
# shifts up by 4 rows.
$ws.Rows("9:12").Delete()

# Type in the new "dealsCategoryTest" block in the now-vacated rows 9-11
# (row 9 retains the yellow-fill header style, row 11 retains the hyperlink
# style, both carried along by the row shift).
$ws.Range("A9").Value = "dealsCategoryTest"
$ws.Range("A10").Value = "Runmode"
$ws.Range("B10").Value = "dealsCategory"
$ws.Range("A11").Value = "Y"
$ws.Range("B11").Value = "Books"

# Move the selection cursor to match the saved workbook state.
$ws.Range("B10").Select()
